$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Code_Only(DeepRL)"
$ws.Range("F9").Select()
